$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark from its original spot
#    (an empty "destinataire" paragraph near the top of the letter,
#    right after the recipient's "Code postal + Ville" placeholder).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Rework the "Je vous informe..." sentence: its paragraph-mark
#    font moves from Arial to Open Sans, and the wording changes.
# ------------------------------------------------------------------
$old = "Je vous informe que copie de ce courrier est transmise " + [char]0x00E0 + " l" + [char]0x2019 + "Inspection du Travail. "
$new = "Je vous informe, que copie de ce courrier est transmise " + [char]0x00E0 + " l" + [char]0x2019 + "inspection du travail, " + [char]0x00E0 + " qui je sollicite, par ailleurs l" + [char]0x2019 + "intervention dans ce dossier"

$searchRange = $d.Content
$found = $searchRange.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $sentencePara = $searchRange.Paragraphs(1)
    $sentencePara.Range.Font.NameAscii = "Open Sans"
    $sentencePara.Range.Font.Name = "Open Sans"
}

$null = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# ------------------------------------------------------------------
# 3) Re-add the "_GoBack" bookmark at the start of the first
#    "« Signature »" paragraph, which immediately follows the
#    sign-off ("Veuillez agréer ... distinguée.") and the employee's
#    name placeholder.
# ------------------------------------------------------------------
$afterSignoff = $d.Content
$null = $afterSignoff.Find.Execute("distingu" + [char]0x00E9 + "e.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$tail = $d.Range($afterSignoff.End, $d.Content.End)
$null = $tail.Find.Execute("Signature", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sigPara = $tail.Paragraphs(1)

$insertPoint = $d.Range($sigPara.Range.Start, $sigPara.Range.Start)
$null = $d.Bookmarks.Add("_GoBack", $insertPoint)
